# Apply the block-order reshuffle edit described by the diff.
# Only the cells that actually change value are touched; everything else
# (styles, untouched cells, dimensions) is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - column labels were reordered.
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "living_rooms_1"
$ws.Range("F1").Value = "living_rooms_2"

# Data rows (the 0/1 indicator matrix) - update only the cells whose value changed.
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = 0

$ws.Range("C4").Value = 0
$ws.Range("F4").Value = 1

$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1

$ws.Range("B6").Value = 1
$ws.Range("F6").Value = 0

$ws.Range("B7").Value = 0
$ws.Range("D7").Value = 1
